$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29; existing rows 29+ shift down by one
# (old row29 -> row30, old row30 -> row31), matching Excel's native
# "insert row" behavior which also carries formatting down correctly.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new data record
$ws.Cells.Item(29, 1).Value = 5
$ws.Cells.Item(29, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(29, 3).Value = "Maule"
$ws.Cells.Item(29, 4).Value = 44474
$ws.Cells.Item(29, 5).Value = 7
$ws.Cells.Item(29, 6).Value = 100112026
$ws.Cells.Item(29, 7).Value = "Haba"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 500
$ws.Cells.Item(29, 11).Value = 8500
$ws.Cells.Item(29, 12).Value = 8500
$ws.Cells.Item(29, 13).Value = 8500
$ws.Cells.Item(29, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(29, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(29, 16).Value = 340
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = "Hortaliza"
